$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" / row 2 (the old first data row). This shifts the
# "Resolving-Mac" row up to become row 2, and the now-unused "ECs" shared
# string will be dropped automatically when the workbook is saved.
$ws.Rows.Item(2).Delete()

# Update row 2 (previously row 3) with the new TPM-derived values.
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.23247666666667
$ws.Range("H2").Value = 60.69743
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5468886666666667
$ws.Range("N2").Value = 1.640666
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 11.06491218759778
$ws.Range("R2").Value = 99.58420968838
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
